# texts.xlsx: "visibility, temperature inside, optimized pngs"
#
# On the "Translation" sheet, the old row 7 (SingleUseId4 / Temperature_wildcard /
# Left / wczytywanie....) is removed, rows 8-10 shift up to 7-9, and a brand new
# row 10 is appended for the new "temperature inside" label:
#   SingleUseId9 | owm_style | Center | Temp. wewnątrz | LTR

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Remove the old row 7 - everything below shifts up by one.
$ws.Rows.Item(7).Delete()

# Append the new translation row for the "temperature inside" text at row 10.
$ws.Range("B10").Value2 = "SingleUseId9"
$ws.Range("C10").Value2 = "owm_style"
$ws.Range("D10").Value2 = "Center"
$ws.Range("E10").Value2 = "Temp. wewnątrz"
$ws.Range("F10").Value2 = "LTR"

# Keep the new row's formatting consistent with the rest of the table (no
# inherited column default style bleeding onto the freshly written cells).
$ws.Range("B10:F10").Style = "Normal"
